{"js": "// Update the worksheet date header and the 25 \"two-digit x two-digit\" answer\n// cells (5 non-empty rows x 5 columns) in the single table to the new values.\n// Replacements are positional (row/col index), NOT global find/replace-by-\n// old-text, because some new values duplicate other old/new values\n// elsewhere in the table (e.g. \"27\u00d766=1782\" and \"36\u00d776=2736\" each appear\n// twice across the before/after sets), so a naive search-and-replace-all\n// would corrupt unrelated cells.\n\n// 1) Date heading paragraph (first paragraph in the document).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2025-02-05 Wednesday\", \"Replace\");\n\n// 2) Table answers: only 5 of the 20 table rows actually hold text\n//    (0-based rows 0, 4, 9, 14, 19); the rest are blank spacer rows and\n//    must stay untouched.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst answers = {\n  0:  [\"70\u00d757=3990\", \"13\u00d786=1118\", \"68\u00d719=1292\", \"45\u00d759=2655\", \"44\u00d767=2948\"],\n  4:  [\"82\u00d769=5658\", \"87\u00d729=2523\", \"36\u00d776=2736\", \"54\u00d716=864\",  \"36\u00d776=2736\"],\n  9:  [\"98\u00d797=9506\", \"99\u00d740=3960\", \"35\u00d750=1750\", \"90\u00d746=4140\", \"82\u00d762=5084\"],\n  14: [\"30\u00d764=1920\", \"47\u00d794=4418\", \"11\u00d799=1089\", \"90\u00d741=3690\", \"43\u00d750=2150\"],\n  19: [\"27\u00d766=1782\", \"50\u00d716=800\",  \"32\u00d785=2720\", \"36\u00d738=1368\", \"35\u00d760=2100\"],\n};\n\nfor (const rowIndex of Object.keys(answers)) {\n  const r = Number(rowIndex);\n  const vals = answers[rowIndex];\n  for (let c = 0; c < vals.length; c++) {\n    table.getCell(r, c).value = vals[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date header and the 25 \"two-digit x two-digit\" answer\n# cells (5 non-empty rows x 5 columns) in the single table to the new values.\n# Replacements are positional (row/col or paragraph index), NOT global\n# find/replace-by-old-text, because some new values duplicate other\n# old/new values elsewhere in the table (e.g. \"27\u00d766=1782\" and\n# \"36\u00d776=2736\" each appear twice across the before/after sets).\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph (first paragraph in the document).\n$d.Paragraphs.Item(1).Range.Text = \"2025-02-05 Wednesday\"\n\n# 2) Table answers: only 5 of the 20 table rows actually hold text\n#    (rows 1, 5, 10, 15, 20 in 1-based Word COM numbering); the rest are\n#    blank spacer rows and must stay untouched.\n$t = $d.Tables.Item(1)\n\n$answers = @{\n  1  = @(\"70\u00d757=3990\", \"13\u00d786=1118\", \"68\u00d719=1292\", \"45\u00d759=2655\", \"44\u00d767=2948\")\n  5  = @(\"82\u00d769=5658\", \"87\u00d729=2523\", \"36\u00d776=2736\", \"54\u00d716=864\",  \"36\u00d776=2736\")\n  10 = @(\"98\u00d797=9506\", \"99\u00d740=3960\", \"35\u00d750=1750\", \"90\u00d746=4140\", \"82\u00d762=5084\")\n  15 = @(\"30\u00d764=1920\", \"47\u00d794=4418\", \"11\u00d799=1089\", \"90\u00d741=3690\", \"43\u00d750=2150\")\n  20 = @(\"27\u00d766=1782\", \"50\u00d716=800\",  \"32\u00d785=2720\", \"36\u00d738=1368\", \"35\u00d760=2100\")\n}\n\nforeach ($r in $answers.Keys) {\n  $vals = $answers[$r]\n  for ($c = 1; $c -le 5; $c++) {\n    $t.Cell($r, $c).Range.Text = $vals[$c - 1]\n  }\n}\n"}
